$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17,8).Value = 2931.0625
$ws.Cells.Item(17,10).Value = 2931.0625
$ws.Cells.Item(17,12).Value = 8793.1875
$ws.Cells.Item(17,14).Value = -9129.1875

$ws.Cells.Item(32,8).Value = 784.625
$ws.Cells.Item(32,10).Value = 896
$ws.Cells.Item(32,12).Value = 896
$ws.Cells.Item(32,14).Value = -1548

$ws.Cells.Item(41,8).Value = 271.42856
$ws.Cells.Item(41,9).Value = 0
$ws.Cells.Item(41,10).Value = 271.42856
$ws.Cells.Item(41,11).Value = 0
$ws.Cells.Item(41,12).Value = 271.42856
$ws.Cells.Item(41,13).Value = ""
$ws.Cells.Item(41,14).Value = -1151.42856

$ws.Cells.Item(62,8).Value = 6131.9165
$ws.Cells.Item(62,9).Value = 4805.5
$ws.Cells.Item(62,11).Value = 4805.5
$ws.Cells.Item(62,13).Value = -4181.5

$ws.Cells.Item(65,8).Value = 6131.9165
$ws.Cells.Item(65,9).Value = 4805.5
$ws.Cells.Item(65,11).Value = 24027.5
$ws.Cells.Item(65,13).Value = -20907.5

$ws.Cells.Item(98,8).Value = 979.8570999999999
$ws.Cells.Item(98,9).Value = 856.9
$ws.Cells.Item(98,10).Value = 1287.25
$ws.Cells.Item(98,11).Value = 856.9
$ws.Cells.Item(98,12).Value = 1287.25
$ws.Cells.Item(98,13).Value = 641.1
$ws.Cells.Item(98,14).Value = -4283.25

$ws.Cells.Item(106,8).Value = 36003.285
$ws.Cells.Item(106,9).Value = 36003.285
$ws.Cells.Item(106,11).Value = 36003.285
$ws.Cells.Item(106,13).Value = -35372.285

$ws.Cells.Item(113,8).Value = 3208
$ws.Cells.Item(113,9).Value = 3272.182
$ws.Cells.Item(113,10).Value = 3090.3333
$ws.Cells.Item(113,11).Value = 3272.182
$ws.Cells.Item(113,12).Value = 3090.3333
$ws.Cells.Item(113,13).Value = -18.18199999999979
$ws.Cells.Item(113,14).Value = -9598.3333

$ws.Cells.Item(115,8).Value = 5210.5
$ws.Cells.Item(115,9).Value = 5210.5
$ws.Cells.Item(115,11).Value = 15631.5
$ws.Cells.Item(115,13).Value = -14064.5

$ws.Cells.Item(122,8).Value = 979.8570999999999
$ws.Cells.Item(122,9).Value = 856.9
$ws.Cells.Item(122,10).Value = 1287.25
$ws.Cells.Item(122,11).Value = 2570.7
$ws.Cells.Item(122,12).Value = 3861.75
$ws.Cells.Item(122,13).Value = -120.6999999999998
$ws.Cells.Item(122,14).Value = -8761.75

$ws.Cells.Item(129,8).Value = 3760.3
$ws.Cells.Item(129,9).Value = 1771.75
$ws.Cells.Item(129,10).Value = 5086
$ws.Cells.Item(129,11).Value = 5315.25
$ws.Cells.Item(129,12).Value = 15258
$ws.Cells.Item(129,13).Value = -315.25
$ws.Cells.Item(129,14).Value = -25258

$ws.Cells.Item(137,8).Value = 2862.5312
$ws.Cells.Item(137,9).Value = 1688.4
$ws.Cells.Item(137,10).Value = 3898.5293
$ws.Cells.Item(137,11).Value = 5065.200000000001
$ws.Cells.Item(137,12).Value = 11695.5879
$ws.Cells.Item(137,13).Value = -2515.200000000001
$ws.Cells.Item(137,14).Value = -16795.5879

$ws.Cells.Item(138,8).Value = 16623.75
$ws.Cells.Item(138,9).Value = 8998.714
$ws.Cells.Item(138,11).Value = 26996.142
$ws.Cells.Item(138,13).Value = -21856.142

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5,8).Value = 366.33334
$ws.Cells.Item(5,9).Value = 49.5
$ws.Cells.Item(5,11).Value = 49.5
$ws.Cells.Item(5,13).Value = 62.5

$ws.Cells.Item(74,8).Value = 3912.5293
$ws.Cells.Item(74,9).Value = 1981.875
$ws.Cells.Item(74,11).Value = 1981.875
$ws.Cells.Item(74,13).Value = -1107.875

$ws.Cells.Item(77,8).Value = 3912.5293
$ws.Cells.Item(77,9).Value = 1981.875
$ws.Cells.Item(77,11).Value = 9909.375
$ws.Cells.Item(77,13).Value = -5541.375

$ws.Cells.Item(96,8).Value = 0
$ws.Cells.Item(96,10).Value = 0
$ws.Cells.Item(96,12).Value = 0
$ws.Cells.Item(96,14).Value = ""

$ws.Cells.Item(97,8).Value = 339.5
$ws.Cells.Item(97,9).Value = 339.5
$ws.Cells.Item(97,10).Value = 0
$ws.Cells.Item(97,11).Value = 339.5
$ws.Cells.Item(97,12).Value = 0
$ws.Cells.Item(97,13).Value = 156.5
$ws.Cells.Item(97,14).Value = ""

$ws.Cells.Item(110,8).Value = 14799.333
$ws.Cells.Item(110,9).Value = 14799.333
$ws.Cells.Item(110,11).Value = 14799.333
$ws.Cells.Item(110,13).Value = -12754.333

$ws.Cells.Item(122,8).Value = 402963.75
$ws.Cells.Item(122,9).Value = 627082.4
$ws.Cells.Item(122,11).Value = 1881247.2
$ws.Cells.Item(122,13).Value = -1878797.2

$ws.Cells.Item(132,8).Value = 1997.2903
$ws.Cells.Item(132,9).Value = 1639.3572
$ws.Cells.Item(132,11).Value = 4918.071599999999
$ws.Cells.Item(132,13).Value = -2388.071599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4,8).Value = 366.33334
$ws.Cells.Item(4,9).Value = 49.5
$ws.Cells.Item(4,11).Value = 49.5
$ws.Cells.Item(4,13).Value = 65.5

$ws.Cells.Item(94,8).Value = 657.38464
$ws.Cells.Item(94,10).Value = 1900
$ws.Cells.Item(94,12).Value = 1900
$ws.Cells.Item(94,14).Value = -2802

$ws.Cells.Item(107,8).Value = 1664.1904
$ws.Cells.Item(107,10).Value = 2778
$ws.Cells.Item(107,12).Value = 2778
$ws.Cells.Item(107,14).Value = -6618

$ws.Cells.Item(134,8).Value = 2383.8462
$ws.Cells.Item(134,9).Value = 777.44446
$ws.Cells.Item(134,11).Value = 2332.33338
$ws.Cells.Item(134,13).Value = 202.66662

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58,8).Value = 8754.200000000001
$ws.Cells.Item(58,9).Value = 8000
$ws.Cells.Item(58,10).Value = 8942.75
$ws.Cells.Item(58,11).Value = 8000
$ws.Cells.Item(58,12).Value = 8942.75
$ws.Cells.Item(58,13).Value = -7797
$ws.Cells.Item(58,14).Value = -9348.75

$ws.Cells.Item(107,8).Value = 1030.1904
$ws.Cells.Item(107,9).Value = 642.875
$ws.Cells.Item(107,11).Value = 642.875
$ws.Cells.Item(107,13).Value = 1277.125

$ws.Cells.Item(120,8).Value = 40000
$ws.Cells.Item(120,10).Value = 40000
$ws.Cells.Item(120,12).Value = 40000
$ws.Cells.Item(120,14).Value = -47258

$ws.Cells.Item(122,8).Value = 644.6923
$ws.Cells.Item(122,9).Value = 615.0833
$ws.Cells.Item(122,11).Value = 1845.2499
$ws.Cells.Item(122,13).Value = 604.7501

$ws.Cells.Item(134,8).Value = 4629.8
$ws.Cells.Item(134,9).Value = 3791.4167
$ws.Cells.Item(134,11).Value = 11374.2501
$ws.Cells.Item(134,13).Value = -8839.250100000001

$ws.Cells.Item(136,8).Value = 8754.200000000001
$ws.Cells.Item(136,9).Value = 8000
$ws.Cells.Item(136,10).Value = 8942.75
$ws.Cells.Item(136,11).Value = 24000
$ws.Cells.Item(136,12).Value = 26828.25
$ws.Cells.Item(136,13).Value = -21450
$ws.Cells.Item(136,14).Value = -31928.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(120,8).Value = 14831.6

$ws.Cells.Item(122,8).Value = 1000.2
$ws.Cells.Item(122,9).Value = 879.8
$ws.Cells.Item(122,11).Value = 7918.2
$ws.Cells.Item(122,13).Value = -5468.2

$ws.Cells.Item(126,8).Value = 2500
$ws.Cells.Item(126,9).Value = 2500
$ws.Cells.Item(126,11).Value = 7500
$ws.Cells.Item(126,13).Value = -2560

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70,8).Value = 6999.6
$ws.Cells.Item(70,10).Value = 6999.6
$ws.Cells.Item(70,12).Value = 6999.6
$ws.Cells.Item(70,14).Value = -7539.6

$ws.Cells.Item(73,8).Value = 6999.6
$ws.Cells.Item(73,10).Value = 6999.6
$ws.Cells.Item(73,12).Value = 6999.6
$ws.Cells.Item(73,14).Value = -8871.6

$ws.Cells.Item(107,8).Value = 1385.7273
$ws.Cells.Item(107,9).Value = 2658.2
$ws.Cells.Item(107,11).Value = 2658.2
$ws.Cells.Item(107,13).Value = -738.1999999999998

$ws.Cells.Item(122,8).Value = 357970.25
$ws.Cells.Item(122,9).Value = 85703.414
$ws.Cells.Item(122,10).Value = 529928.25
$ws.Cells.Item(122,11).Value = 257110.242
$ws.Cells.Item(122,12).Value = 1589784.75
$ws.Cells.Item(122,13).Value = -254660.242
$ws.Cells.Item(122,14).Value = -1594684.75

$ws.Cells.Item(126,8).Value = 4298.1665
$ws.Cells.Item(126,9).Value = 3184
$ws.Cells.Item(126,10).Value = 4855.25
$ws.Cells.Item(126,11).Value = 9552
$ws.Cells.Item(126,12).Value = 14565.75
$ws.Cells.Item(126,13).Value = -7082
$ws.Cells.Item(126,14).Value = -19505.75

$ws.Cells.Item(132,8).Value = 3073.6365
$ws.Cells.Item(132,9).Value = 2360.2104
$ws.Cells.Item(132,11).Value = 7080.6312
$ws.Cells.Item(132,13).Value = -4550.6312

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7,8).Value = 3981.4
$ws.Cells.Item(7,9).Value = 3981.4
$ws.Cells.Item(7,11).Value = 3981.4
$ws.Cells.Item(7,13).Value = -3869.4

$ws.Cells.Item(22,8).Value = 1998.2
$ws.Cells.Item(22,9).Value = 1998.25
$ws.Cells.Item(22,10).Value = 1998
$ws.Cells.Item(22,11).Value = 1998.25
$ws.Cells.Item(22,12).Value = 1998
$ws.Cells.Item(22,13).Value = -1703.25
$ws.Cells.Item(22,14).Value = -2588

$ws.Cells.Item(27,8).Value = 1998.2
$ws.Cells.Item(27,9).Value = 1998.25
$ws.Cells.Item(27,10).Value = 1998
$ws.Cells.Item(27,11).Value = 1998.25
$ws.Cells.Item(27,12).Value = 1998
$ws.Cells.Item(27,13).Value = -1891.25
$ws.Cells.Item(27,14).Value = -2212

$ws.Cells.Item(40,8).Value = 4751.857
$ws.Cells.Item(40,9).Value = 4727.1665
$ws.Cells.Item(40,11).Value = 4727.1665
$ws.Cells.Item(40,13).Value = -4591.1665

$ws.Cells.Item(109,8).Value = 52210
$ws.Cells.Item(109,10).Value = 52000
$ws.Cells.Item(109,12).Value = 52000
$ws.Cells.Item(109,14).Value = -54774

$ws.Cells.Item(126,8).Value = 3981.4
$ws.Cells.Item(126,9).Value = 3981.4
$ws.Cells.Item(126,11).Value = 11944.2
$ws.Cells.Item(126,13).Value = -9474.200000000001

$ws.Cells.Item(132,8).Value = 4647.3716
$ws.Cells.Item(132,9).Value = 3173.5
$ws.Cells.Item(132,11).Value = 9520.5
$ws.Cells.Item(132,13).Value = -6990.5

$ws.Cells.Item(136,8).Value = 4081.0908
$ws.Cells.Item(136,9).Value = 3989.2
$ws.Cells.Item(136,11).Value = 11967.6
$ws.Cells.Item(136,13).Value = -9417.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96,8).Value = 901
$ws.Cells.Item(96,10).Value = 711.25
$ws.Cells.Item(96,12).Value = 711.25
$ws.Cells.Item(96,14).Value = -3457.25

$ws.Cells.Item(122,8).Value = 2097.6667
$ws.Cells.Item(122,9).Value = 2047.375
$ws.Cells.Item(122,10).Value = 2500
$ws.Cells.Item(122,11).Value = 6142.125
$ws.Cells.Item(122,12).Value = 7500
$ws.Cells.Item(122,13).Value = -3692.125
$ws.Cells.Item(122,14).Value = -12400
